$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the column-G JSON-builder formulas -----------------------------
# Original formula (per row N):
#   CHAR(34)&"PuzzleNum"&CHAR(34)&": "&CHAR(34)&AN&CHAR(34)&","&CHAR(34)&"Title"&CHAR(34)&": "&CHAR(34)&BN&CHAR(34)&","&CHAR(34)&"Pieces"&CHAR(34)&": "&CHAR(34)&CN&CHAR(34)&","&CHAR(34)&"Company"&CHAR(34)&": "&CHAR(34)&DN&CHAR(34)&","&CHAR(34)&"Size"&CHAR(34)&": "&CHAR(34)&EN&CHAR(34)&","&CHAR(34)&"URL"&CHAR(34)&": "&CHAR(34)&FN&CHAR(34)
# New formula wraps that in a JSON object literal and trailing comma:
#   "{"&CHAR(34)&"PuzzleNum"&CHAR(34)&...&CHAR(34)&FN&CHAR(34)&"},"
#
# Row 2 is a standalone (non-shared) formula in the source file, so it is set
# separately from the rest (G3:G176) so the shared-formula grouping on save
# matches the original layout (G3:G66, G67:G130, G131:G176).

$ws.Range("G2").Formula = "=""{""&CHAR(34)&""PuzzleNum""&CHAR(34)&"": ""&CHAR(34)&A2&CHAR(34)&"",""&CHAR(34)&""Title""&CHAR(34)&"": ""&CHAR(34)&B2&CHAR(34)&"",""&CHAR(34)&""Pieces""&CHAR(34)&"": ""&CHAR(34)&C2&CHAR(34)&"",""&CHAR(34)&""Company""&CHAR(34)&"": ""&CHAR(34)&D2&CHAR(34)&"",""&CHAR(34)&""Size""&CHAR(34)&"": ""&CHAR(34)&E2&CHAR(34)&"",""&CHAR(34)&""URL""&CHAR(34)&"": ""&CHAR(34)&F2&CHAR(34)&""},"""

$ws.Range("G3:G176").Formula = "=""{""&CHAR(34)&""PuzzleNum""&CHAR(34)&"": ""&CHAR(34)&A3&CHAR(34)&"",""&CHAR(34)&""Title""&CHAR(34)&"": ""&CHAR(34)&B3&CHAR(34)&"",""&CHAR(34)&""Pieces""&CHAR(34)&"": ""&CHAR(34)&C3&CHAR(34)&"",""&CHAR(34)&""Company""&CHAR(34)&"": ""&CHAR(34)&D3&CHAR(34)&"",""&CHAR(34)&""Size""&CHAR(34)&"": ""&CHAR(34)&E3&CHAR(34)&"",""&CHAR(34)&""URL""&CHAR(34)&"": ""&CHAR(34)&F3&CHAR(34)&""},"""

# --- Update the view / selection -------------------------------------------
# Scroll the window so row 163 is (as close as possible to) the top-left cell
# and move the active selection to G2:G176 with G2 as the active cell.
$excel.Goto($ws.Range("C163"), $true)
$ws.Range("G2:G176").Select()
